$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) hold text values that look numeric
# (e.g. "1.004", "48.28") or percentages. Force text format so Excel does
# not reinterpret them as numbers/dates when we assign new string values.
$ws.Range("D2:E51").NumberFormat = "@"

function Set-Text($cellRef, $text) {
    $ws.Range($cellRef).Value = $text
}

# Row 2 - Bitcoin
Set-Text "D2" "27.452.85"
Set-Text "E2" "  -3.46%  "

# Row 3 - Ethereum
Set-Text "D3" "1.853.04"
Set-Text "E3" "  -5.10%  "

# Row 4 - TetherUSD
Set-Text "D4" "1.003"
Set-Text "E4" "  -0.86%  "

# Row 5 - BNB
Set-Text "D5" "320.91"
Set-Text "E5" "  +0.04%  "

# Row 6 - USDC
Set-Text "D6" "1.003"
Set-Text "E6" "  -0.74%  "

# Row 7 - XRP
Set-Text "D7" "0.4484"
Set-Text "E7" "  -5.73%  "

# Row 8 - Cardano
Set-Text "D8" "0.3846"
Set-Text "E8" "  -5.44%  "

# Row 9 - OKB
Set-Text "D9" "48.30"
Set-Text "E9" "  -9.61%  "

# Row 10 - Dogecoin
Set-Text "D10" "0.07840"
Set-Text "E10" "  -7.30%  "

# Row 11 - Polygon
Set-Text "E11" "  -3.76%  "

# Row 12 - Solana
Set-Text "E12" "  -2.95%  "

# Row 13 - was WrappedEther, now Polkadot
Set-Text "B13" "Polkadot"
Set-Text "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-Text "D13" "5.860"
Set-Text "E13" "  -4.95%  "

# Row 14 - was Polkadot, now Chainlink
Set-Text "B14" "Chainlink"
Set-Text "C14" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-Text "D14" "7.107"
Set-Text "E14" "  -6.51%  "

# Row 15 - was Chainlink, now WrappedEther
Set-Text "B15" "WrappedEther"
Set-Text "C15" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-Text "D15" "1.777.51"
Set-Text "E15" "  -10.27%  "

# Row 16 - BinanceUSD
Set-Text "E16" "  -0.72%  "

# Row 17 - Litecoin
Set-Text "D17" "85.96"
Set-Text "E17" "  -3.52%  "

# Row 18 - ShibaInu
Set-Text "D18" "0.00001024"
Set-Text "E18" "  -4.47%  "

# Row 19 - TRON
Set-Text "D19" "0.06517"
Set-Text "E19" "  -1.55%  "

# Row 20 - Avalanche
Set-Text "D20" "17.08"
Set-Text "E20" "  -8.68%  "

# Row 21 - Dai
Set-Text "D21" "1.003"
Set-Text "E21" "  -0.86%  "

# Row 22 - Uniswap
Set-Text "D22" "5.490"
Set-Text "E22" "  -5.70%  "

# Row 23 - WrappedBTC
Set-Text "D23" "27.409.58"
Set-Text "E23" "  -3.64%  "

# Row 24 - Cosmos
Set-Text "D24" "10.79"
Set-Text "E24" "  -6.38%  "

# Row 25 - Toncoin
Set-Text "D25" "2.295"
Set-Text "E25" "  +0.17%  "

# Row 26 - WrappedliquidstakedEther2.0
Set-Text "D26" "2.043.07"
Set-Text "E26" "  -7.70%  "

# Row 27 - Monero
Set-Text "D27" "150.65"
Set-Text "E27" "  -2.15%  "

# Row 28 - EthereumClassic
Set-Text "E28" "  -4.26%  "

# Row 29 - InternetComputer(DFINITY)
Set-Text "D29" "5.479"
Set-Text "E29" "  -8.09%  "

# Row 30 - LidoDAOToken
Set-Text "E30" "  -5.92%  "

# Row 31 - BitcoinCash
Set-Text "D31" "119.86"
Set-Text "E31" "  -3.08%  "

# Row 32 - ARBITRUM
Set-Text "D32" "1.499"
Set-Text "E32" "  +3.49%  "

# Row 33 - Stellar
Set-Text "D33" "0.09360"
Set-Text "E33" "  -2.26%  "

# Row 34 - ImmutableX
Set-Text "D34" "0.9261"
Set-Text "E34" "  -6.13%  "

# Row 35 - HuobiToken
Set-Text "D35" "3.615"
Set-Text "E35" "  -0.94%  "

# Row 36 - Filecoin
Set-Text "D36" "5.239"
Set-Text "E36" "  -6.37%  "

# Row 37 - TrustWalletToken
Set-Text "D37" "1.223"
Set-Text "E37" "  -2.54%  "

# Row 38 - VeChain
Set-Text "D38" "0.02217"
Set-Text "E38" "  -4.97%  "

# Row 39 - Hedera
Set-Text "D39" "0.05939"
Set-Text "E39" "  -4.37%  "

# Row 40 - FraxShare
Set-Text "D40" "8.292"
Set-Text "E40" "  -5.81%  "

# Row 41 - Frax
Set-Text "D41" "1.002"

# Row 42 - TheSandbox
Set-Text "D42" "0.5898"
Set-Text "E42" "  -5.18%  "

# Row 43 - Aptos
Set-Text "D43" "10.28"
Set-Text "E43" "  -7.60%  "

# Row 44 - Algorand
Set-Text "D44" "0.1843"
Set-Text "E44" "  -3.99%  "

# Row 45 - WEMIXTOKEN
Set-Text "D45" "1.279"
Set-Text "E45" "  -4.52%  "

# Row 46 - Decentraland
Set-Text "D46" "0.5615"
Set-Text "E46" "  -5.80%  "

# Row 47 - EnergySwap
Set-Text "D47" "12.19"
Set-Text "E47" "  -6.38%  "

# Row 48 - PancakeSwap
Set-Text "D48" "3.353"

# Row 49 - NEARProtocol
Set-Text "D49" "1.912"
Set-Text "E49" "  -6.81%  "

# Row 50 - Cronos
Set-Text "D50" "0.06839"
Set-Text "E50" "  +0.29%  "

# Row 51 - PaxosStandard
Set-Text "D51" "1.003"
Set-Text "E51" "  -11.04%  "
